$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.498.43"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.34"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.17"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4558"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3661"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07126"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8788"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.35"
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.810.72"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.286"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.372"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.54"
$ws.Range("E16").Value = "  -5.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008577"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.563.02"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.009"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.42"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.93"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.059"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.74"
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.837"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08678"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.506"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7312"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.118"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.661"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01950"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05106"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.897"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.965"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4985"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1559"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.160"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4597"
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.941"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.92"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.589"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06000"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.34"
$ws.Range("E51").Value = "  -1.54%  "
